$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G:H (shifts old "fantasy points" column from G to I)
$ws.Range("G1:H1").EntireColumn.Insert()

# New headers for the inserted columns
$ws.Range("G1").Value = "height"
$ws.Range("H1").Value = "weight"

# Fill height/weight values for each data row (2 through 17)
$lastRow = 17
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 6.166666666666667
    $ws.Cells.Item($r, 8).Value = 220
}
